# "Lagt til CCS og oppdatert scengen" - add CCS-related rows and update scenarios
$wb = $excel.ActiveWorkbook

# Power systems: extend from 12 values (A1:A12) to 14 values (A1:A14)
$wsPower = $wb.Worksheets.Item("Power systems")
$wsPower.Cells.Item(13, 1).Value = 12
$wsPower.Cells.Item(14, 1).Value = 13

# Scenarios: extend from 1 value (A1:A2) to 27 values (A1:A28)
$wsScenarios = $wb.Worksheets.Item("Scenarios")
For ($i = 3; $i -le 28; $i++) {
    $wsScenarios.Cells.Item($i, 1).Value = $i - 1
}

# Routes: shrink from 1062 values (A1:A1063) to 858 values (A1:A859)
$wsRoutes = $wb.Worksheets.Item("Routes")
$wsRoutes.Range("A860:A1063").Clear()

# Installations: shrink from 16 values (A1:A17) to 11 values (A1:A12)
$wsInstallations = $wb.Worksheets.Item("Installations")
$wsInstallations.Range("A13:A17").Clear()
